# LDLC smartphone price tracker: a new hourly price-snapshot column is
# appended to the price-history block. The sheet stores one column per
# scrape timestamp (B..FR), followed by "nom" (product name) and
# "url_produit" (product url). A new snapshot column is inserted right
# before "nom", pushing "nom" -> FT and "url_produit" -> FU, and its
# values are the latest known price for each product (a copy of the
# previous last snapshot column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column FS (175) is where "nom" currently lives; column FR (174) is the
# most recent existing price-snapshot column. Inserting a whole column at
# FS shifts "nom" -> FT and "url_produit" -> FU, and shifts nothing to the
# left of FS (so FR keeps holding the previous last snapshot's values).
$lastSnapshotCol = 174
$newSnapshotCol = 175

$ws.Range("FS1").EntireColumn.Insert()

# Header for the freshly inserted column: the new snapshot's timestamp.
$ws.Range("FS1").Value = "2026-02-05 01:01:02"

# Every data row: carry the latest known price forward into the new
# snapshot column (rows with no price yet simply copy the blank).
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $lastVal = $ws.Cells.Item($r, $lastSnapshotCol).Value2
    $ws.Cells.Item($r, $newSnapshotCol).Value = $lastVal
}
